# Slide 2 ("Agenda") contains three shapes that make up the right-hand
# "EXPLORING RESEARCH QUESTIONS" card: a rounded-rectangle card
# background, the network-diagram icon picture inside it, and the
# "SYSTEM GRAPHIC" caption textbox. All three shift left by 55560 EMU
# (~4.37 pt) to re-center the card's contents; nothing else changes.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Shape 5: "Rounded Rectangle 20" card background -> new x = 5207891 EMU
$card = $s.Shapes.Item(5)
$card.Left = 410.07015748031495

# Shape 6: "Graphic 19" network-diagram icon picture -> new x = 5583240 EMU
$icon = $s.Shapes.Item(6)
$icon.Left = 439.6251988503937

# Shape 7: "TextBox 21" - "SYSTEM GRAPHIC" caption -> new x = 5217309 EMU
$caption = $s.Shapes.Item(7)
$caption.Left = 410.81173228346455
